$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2025-07-31 Thursday" "2025-08-01 Friday"

Replace-Text "58×70=" "29×79="
Replace-Text "82×61=" "21×30="
Replace-Text "30×77=" "74×71="
Replace-Text "93×69=" "80×73="
Replace-Text "87×39=" "11×22="
Replace-Text "38×88=" "59×18="
Replace-Text "50×98=" "45×84="
Replace-Text "67×38=" "27×25="
Replace-Text "44×42=" "13×96="
Replace-Text "37×25=" "44×17="
Replace-Text "27×97=" "38×76="
Replace-Text "55×58=" "43×91="
Replace-Text "38×26=" "21×14="
Replace-Text "76×37=" "54×97="
Replace-Text "20×27=" "73×75="
Replace-Text "88×38=" "59×87="
Replace-Text "18×59=" "20×67="
Replace-Text "34×95=" "61×67="
Replace-Text "43×57=" "37×22="
Replace-Text "82×47=" "26×36="
Replace-Text "76×85=" "39×89="
Replace-Text "71×72=" "16×87="
Replace-Text "45×98=" "27×90="
Replace-Text "99×86=" "59×28="
Replace-Text "76×94=" "84×70="
